$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.65000000000057
$ws.Range("H2").Value = [double]"2.846725704167068e-16"
$ws.Range("K2").Value = 53.40149114864229
$ws.Range("L2").Value = "[44.34173254075646, 62.46124975652811]"
$ws.Range("O2").Value = 1.628973968528041
$ws.Range("P2").Value = "[1.4402897250691948, 1.8176582119868874]"
$ws.Range("S2").Value = 64.78791272394777
$ws.Range("T2").Value = "[58.93545768119077, 70.64036776670477]"
$ws.Range("W2").Value = 19.00000000000042
$ws.Range("X2").Value = 18.22972972973013
$ws.Range("Y2").Value = 19.77027027027071

# Row 3
$ws.Range("E3").Value = 25.8300000000006
$ws.Range("H3").Value = [double]"2.846725704167068e-16"
$ws.Range("K3").Value = 44.04135933414456
$ws.Range("L3").Value = "[33.332781496590876, 54.74993717169824]"
$ws.Range("M3").Value = [double]"1.110223024625157e-14"
$ws.Range("N3").Value = [double]"1.110223024625157e-14"
$ws.Range("O3").Value = 2.132131951084965
$ws.Range("P3").Value = "[1.8931319093704264, 2.371131992799503]"
$ws.Range("S3").Value = 61.22989757301215
$ws.Range("T3").Value = "[55.572140454835306, 66.88765469118898]"
$ws.Range("W3").Value = 17.06486486486526
$ws.Range("X3").Value = 16.08234234234272
$ws.Range("Y3").Value = 18.0473873873878
